# animation refactor (i love refactoring for the sake of readability and easiness)
#
# - remove the empty leftover "Tabelle2" sheet, leaving "Tabelle1" as the
#   only (and active) worksheet
# - log a new time-tracking entry (row 17) for Aris
# - update the saved view state (zoom + selection) to match where work left off

$wb = $excel.ActiveWorkbook

# Suppress the "this sheet contains data, delete anyway?" style prompt.
$excel.DisplayAlerts = $false

$tabelle2 = $wb.Worksheets.Item("Tabelle2")
$tabelle2.Delete()

$ws = $wb.Worksheets.Item("Tabelle1")

# --- new time log entry: Aris, 18:30 -> 19:42, GitHub commit about the animation refactor ---
$ws.Range("A17").Value = "Aris"

$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)            # xlPasteFormats - reuse the date number format
$ws.Range("B17").Value = 45326

$ws.Range("C17").Formula = "= 18+30/60"
$ws.Range("D17").Formula = "=19+42/60"

$ws.Range("F17").Value = "refactor animation -> slip jumping and falling + merged playerMovement and animController"

# --- restore the workbook/sheet view the author ended up with ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("H19").Select()

$wb.Save()
